# Allowing exception file from EHS to be processed (#949)
# data file is put in ENV['CAMPUS_ACCESS_DIRECTORY']

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page1_1")

# Decrement the three placeholder SSN values (999999999) for rows 6, 8 and 9
# so each row gets its own unique exception identifier.
$ws.Range("H6").Value = 999999998
$ws.Range("H8").Value = 999999997
$ws.Range("H9").Value = 999999996

# Give the real SSN value in H7 a dedicated 9-digit zero-padded display
# format (000000000), matching the new number format / style added to the
# workbook.
$ws.Range("H7").NumberFormat = "000000000"

# Update the active selection to reflect where the edit was made.
$ws.Range("H8").Select() | Out-Null
